$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sigma_010")
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 28.10792215963118
$ws.Cells.Item(2, 3).Value = 32.05345581224955
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 28.11777440470371
$ws.Cells.Item(3, 3).Value = 32.09127818402621
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 28.14381941230749
$ws.Cells.Item(4, 3).Value = 32.05462947410271
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 28.12480693753092
$ws.Cells.Item(5, 3).Value = 32.0301315161954
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 28.17516414984882
$ws.Cells.Item(6, 3).Value = 32.07160624309774
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 28.1531843665557
$ws.Cells.Item(7, 3).Value = 32.05443232977905
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 28.14640895792944
$ws.Cells.Item(8, 3).Value = 32.05030447396103
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 28.15550935647234
$ws.Cells.Item(9, 3).Value = 32.05121752348959
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 28.11680450606006
$ws.Cells.Item(10, 3).Value = 32.04951812926105
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 28.15215382069112
$ws.Cells.Item(11, 3).Value = 32.06297356399036
$ws.Cells.Item(12, 2).Value = 28.13935480717308
$ws.Cells.Item(12, 3).Value = 32.05695472501527

$ws = $wb.Worksheets.Item("sigma_025")
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 19.728898764027
$ws.Cells.Item(2, 3).Value = 28.15347852832511
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 19.70595960837813
$ws.Cells.Item(3, 3).Value = 28.13703451901332
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 19.72178555133004
$ws.Cells.Item(4, 3).Value = 28.13403374384782
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 19.72064166612591
$ws.Cells.Item(5, 3).Value = 28.16424832429117
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 19.71297852672667
$ws.Cells.Item(6, 3).Value = 28.14152745576934
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 19.74178674772681
$ws.Cells.Item(7, 3).Value = 28.16674317199548
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 19.73285418503712
$ws.Cells.Item(8, 3).Value = 28.13663119091968
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 19.72752058853814
$ws.Cells.Item(9, 3).Value = 28.15672244995626
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 19.72224611023843
$ws.Cells.Item(10, 3).Value = 28.13562197306298
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 19.69238771847423
$ws.Cells.Item(11, 3).Value = 28.10207817607494
$ws.Cells.Item(12, 2).Value = 19.72070594666025
$ws.Cells.Item(12, 3).Value = 28.14281195332561

$ws = $wb.Worksheets.Item("sigma_050")
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 14.70040253248522
$ws.Cells.Item(2, 3).Value = 23.69408789450418
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 14.67537168995946
$ws.Cells.Item(3, 3).Value = 23.67196706121245
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 14.69562797425608
$ws.Cells.Item(4, 3).Value = 23.67607795289606
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 14.67355480684625
$ws.Cells.Item(5, 3).Value = 23.65305598500715
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 14.70840102243224
$ws.Cells.Item(6, 3).Value = 23.73993404927279
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 14.67130193969093
$ws.Cells.Item(7, 3).Value = 23.71905555324877
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 14.66907123849355
$ws.Cells.Item(8, 3).Value = 23.62678657864157
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 14.67519272858435
$ws.Cells.Item(9, 3).Value = 23.67257957510466
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 14.68538070810352
$ws.Cells.Item(10, 3).Value = 23.7121140070005
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 14.68009263299783
$ws.Cells.Item(11, 3).Value = 23.69412438775136
$ws.Cells.Item(12, 2).Value = 14.68343972738494
$ws.Cells.Item(12, 3).Value = 23.68597830446395
